$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 17
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 18
$ws.Range("B9").Value = 8
$ws.Range("B10").Value = 15
$ws.Range("B14").Value = 8
$ws.Range("B15").Value = 8
$ws.Range("B16").Value = 7
$ws.Range("B17").Value = 9
$ws.Range("B19").Value = 16
$ws.Range("B23").Value = 8
$ws.Range("B25").Value = 14
$ws.Range("B26").Value = 8
